# Apply value updates to the Typhon_Profits workbook (Sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each affected leve row, set the new currentAveragePrice / LevePrice* / LeveProfit* figures,
# clearing cells that no longer hold a value and adding cells that are newly populated.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1534.4
$ws.Range("I40").Value = 1137
$ws.Range("J40").Value = 2130.5
$ws.Range("K40").Value = 1137
$ws.Range("L40").Value = 2130.5
$ws.Range("M40").Value = -962
$ws.Range("N40").Value = -2480.5

$ws.Range("H43").Value = 261.75
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 261.75
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 261.75
$ws.Range("N43").Value = -399.75
$ws.Range("M43").ClearContents()

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H116").Value = 3940.3125
$ws.Range("I116").Value = 2129.875
$ws.Range("K116").Value = 2129.875
$ws.Range("M116").Value = 1312.125

$ws.Range("H135").Value = 26324608
$ws.Range("I135").Value = 896
$ws.Range("J135").Value = 100031000
$ws.Range("K135").Value = 8064
$ws.Range("L135").Value = 900279000
$ws.Range("M135").Value = -5529
$ws.Range("N135").Value = -900284070

$ws.Range("H137").Value = 1531.4138
$ws.Range("I137").Value = 1362.9584
$ws.Range("J137").Value = 2340
$ws.Range("K137").Value = 4088.8752
$ws.Range("L137").Value = 7020
$ws.Range("M137").Value = -1538.8752
$ws.Range("N137").Value = -12120


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1449.3414
$ws.Range("I2").Value = 1332.0286
$ws.Range("K2").Value = 1332.0286
$ws.Range("M2").Value = -1219.0286

$ws.Range("H32").Value = 6650.7393
$ws.Range("I32").Value = 5120.3657
$ws.Range("J32").Value = 19199.8
$ws.Range("K32").Value = 5120.3657
$ws.Range("L32").Value = 19199.8
$ws.Range("M32").Value = -4833.3657
$ws.Range("N32").Value = -19773.8

$ws.Range("H61").Value = 1922.9642
$ws.Range("I61").Value = 1699
$ws.Range("J61").Value = 2953.2
$ws.Range("K61").Value = 1699
$ws.Range("L61").Value = 2953.2
$ws.Range("M61").Value = -1487
$ws.Range("N61").Value = -3377.2

$ws.Range("H116").Value = 1449.3414
$ws.Range("I116").Value = 1332.0286
$ws.Range("K116").Value = 1332.0286
$ws.Range("M116").Value = 961.9713999999999

$ws.Range("H132").Value = 13798.269
$ws.Range("I132").Value = 1452.1892
$ws.Range("K132").Value = 4356.5676
$ws.Range("M132").Value = -1826.5676

$ws.Range("H136").Value = 1922.9642
$ws.Range("I136").Value = 1699
$ws.Range("J136").Value = 2953.2
$ws.Range("K136").Value = 5097
$ws.Range("L136").Value = 8859.599999999999
$ws.Range("M136").Value = -2547
$ws.Range("N136").Value = -13959.6


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1449.3414
$ws.Range("I3").Value = 1332.0286
$ws.Range("K3").Value = 1332.0286
$ws.Range("M3").Value = -1218.0286

$ws.Range("H11").Value = 287
$ws.Range("I11").Value = 287
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 287
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -147
$ws.Range("N11").ClearContents()

$ws.Range("H20").Value = 5011.5
$ws.Range("I20").Value = 6332
$ws.Range("J20").Value = 1050
$ws.Range("K20").Value = 6332
$ws.Range("L20").Value = 1050
$ws.Range("M20").Value = -6085
$ws.Range("N20").Value = -1544

$ws.Range("H76").Value = 25000
$ws.Range("J76").Value = 25000
$ws.Range("L76").Value = 25000
$ws.Range("N76").Value = -25630

$ws.Range("H79").Value = 25000
$ws.Range("J79").Value = 25000
$ws.Range("L79").Value = 25000
$ws.Range("N79").Value = -27184

$ws.Range("H134").Value = 5289.6523
$ws.Range("I134").Value = 6025.3687
$ws.Range("J134").Value = 1795
$ws.Range("K134").Value = 18076.1061
$ws.Range("L134").Value = 5385
$ws.Range("M134").Value = -15541.1061
$ws.Range("N134").Value = -10455


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17509.291
$ws.Range("J31").Value = 5123.4443
$ws.Range("L31").Value = 5123.4443
$ws.Range("N31").Value = -5713.4443

$ws.Range("H34").Value = 17509.291
$ws.Range("J34").Value = 5123.4443
$ws.Range("L34").Value = 5123.4443
$ws.Range("N34").Value = -5527.4443

$ws.Range("H132").Value = 13180.674
$ws.Range("I132").Value = 17427.562
$ws.Range("J132").Value = 3473.5
$ws.Range("K132").Value = 52282.686
$ws.Range("L132").Value = 10420.5
$ws.Range("M132").Value = -49752.686
$ws.Range("N132").Value = -15480.5

$ws.Range("H134").Value = 860.8421
$ws.Range("I134").Value = 736.13336
$ws.Range("K134").Value = 2208.40008
$ws.Range("M134").Value = 326.5999199999997


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1002.5
$ws.Range("I5").Value = 1077.5
$ws.Range("J5").Value = 852.5
$ws.Range("K5").Value = 3232.5
$ws.Range("L5").Value = 2557.5
$ws.Range("M5").Value = -3120.5
$ws.Range("N5").Value = -2781.5

$ws.Range("H34").Value = 928.9286
$ws.Range("J34").Value = 988.5833
$ws.Range("L34").Value = 2965.7499
$ws.Range("N34").Value = -3133.7499

$ws.Range("H37").Value = 35770572
$ws.Range("J37").Value = 35770572
$ws.Range("L37").Value = 107311716
$ws.Range("N37").Value = -107311940

$ws.Range("H39").Value = 2950
$ws.Range("J39").Value = 2950
$ws.Range("L39").Value = 8850
$ws.Range("N39").Value = -9438

$ws.Range("H51").Value = 2245
$ws.Range("I51").Value = 2000
$ws.Range("J51").Value = 2490
$ws.Range("K51").Value = 6000
$ws.Range("L51").Value = 7470
$ws.Range("M51").Value = -5540
$ws.Range("N51").Value = -8390

$ws.Range("H55").Value = 3300
$ws.Range("J55").Value = 3300
$ws.Range("L55").Value = 9900
$ws.Range("N55").Value = -10254

$ws.Range("H87").Value = 12987.533
$ws.Range("I87").Value = 7092.091
$ws.Range("J87").Value = 29200
$ws.Range("K87").Value = 21276.273
$ws.Range("L87").Value = 87600
$ws.Range("M87").Value = -20028.273
$ws.Range("N87").Value = -90096

$ws.Range("H90").Value = 12987.533
$ws.Range("I90").Value = 7092.091
$ws.Range("J90").Value = 29200
$ws.Range("K90").Value = 63828.819
$ws.Range("L90").Value = 262800
$ws.Range("M90").Value = -57588.819
$ws.Range("N90").Value = -275280

$ws.Range("H107").Value = 4256
$ws.Range("I107").Value = 20158
$ws.Range("K107").Value = 60474
$ws.Range("M107").Value = -58554

$ws.Range("H110").Value = 5150
$ws.Range("I110").Value = 300
$ws.Range("K110").Value = 900
$ws.Range("M110").Value = 3190

$ws.Range("H113").Value = 9712.546
$ws.Range("J113").Value = 639.5714
$ws.Range("L113").Value = 1918.7142
$ws.Range("N113").Value = -6258.7142

$ws.Range("H131").Value = 760.7
$ws.Range("I131").Value = 283
$ws.Range("J131").Value = 775.47424
$ws.Range("K131").Value = 849
$ws.Range("L131").Value = 2326.42272
$ws.Range("M131").Value = 4191
$ws.Range("N131").Value = -12406.42272

$ws.Range("H135").Value = 1002.5
$ws.Range("I135").Value = 1077.5
$ws.Range("J135").Value = 852.5
$ws.Range("K135").Value = 9697.5
$ws.Range("L135").Value = 7672.5
$ws.Range("M135").Value = -7162.5
$ws.Range("N135").Value = -12742.5

$ws.Range("H139").Value = 2255.0588
$ws.Range("I139").Value = 1258.0358
$ws.Range("K139").Value = 3774.1074
$ws.Range("M139").Value = 1365.8926


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 21900
$ws.Range("J92").Value = 21900
$ws.Range("L92").Value = 21900
$ws.Range("N92").Value = -25644

$ws.Range("H126").Value = 5183.3335
$ws.Range("I126").Value = 4083.3333
$ws.Range("K126").Value = 12249.9999
$ws.Range("M126").Value = -9779.999899999999

$ws.Range("H132").Value = 18260.47
$ws.Range("I132").Value = 3698.4443
$ws.Range("J132").Value = 74428.28999999999
$ws.Range("K132").Value = 11095.3329
$ws.Range("L132").Value = 223284.87
$ws.Range("M132").Value = -8565.332900000001
$ws.Range("N132").Value = -228344.87


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 687945.9399999999
$ws.Range("I2").Value = 854999.7
$ws.Range("J2").Value = 131100
$ws.Range("K2").Value = 854999.7
$ws.Range("L2").Value = 131100
$ws.Range("M2").Value = -854887.7
$ws.Range("N2").Value = -131324

$ws.Range("H68").Value = 1212.6364
$ws.Range("I68").Value = 1212.6364
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1212.6364
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -463.6364000000001
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1212.6364
$ws.Range("I71").Value = 1212.6364
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 6063.182000000001
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -2319.182000000001
$ws.Range("N71").ClearContents()

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H132").Value = 1749.6471
$ws.Range("I132").Value = 1218.6666
$ws.Range("J132").Value = 2607.3845
$ws.Range("K132").Value = 3655.9998
$ws.Range("L132").Value = 7822.1535
$ws.Range("M132").Value = -1125.9998
$ws.Range("N132").Value = -12882.1535


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 11779.2
$ws.Range("J69").Value = 13723.75
$ws.Range("L69").Value = 13723.75
$ws.Range("N69").Value = -15221.75

$ws.Range("H72").Value = 11779.2
$ws.Range("J72").Value = 13723.75
$ws.Range("L72").Value = 41171.25
$ws.Range("N72").Value = -48659.25

$ws.Range("H81").Value = 111112880
$ws.Range("I81").Value = 1991.375
$ws.Range("K81").Value = 3982.75
$ws.Range("M81").Value = -2921.75

$ws.Range("H84").Value = 111112880
$ws.Range("I84").Value = 1991.375
$ws.Range("K84").Value = 19913.75
$ws.Range("M84").Value = -14609.75

